$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9743480086326599
$ws.Range("B1").Value = 0.9325898885726929
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.587123155593872
$ws.Range("E1").Value = 0.9593753814697266
